# DSBot/kb.xlsx - add regression / cross validation (kFold) / tweak classification & roc curve selections
$wb = $excel.ActiveWorkbook

# --- Foglio2: just move the current selection (no content change) ---
$ws2 = $wb.Worksheets.Item("Foglio2")
$ws2.Range("A123:G128").Select()

# --- Foglio4: move the current selection (no content change) ---
$ws3 = $wb.Worksheets.Item("Foglio4")
$ws3.Range("A214").Select()

# --- Foglio3: move the current selection (no content change) ---
$ws4 = $wb.Worksheets.Item("Foglio3")
$ws4.Range("A5").Select()

# --- New sheet Foglio5: cross validation (kFold) step added to the
#     missingValues/zeroVariance/categorical/hasLabel/moreFeatures/outliers
#     combinations, right before the randomForest + roc pair ---
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "Foglio5"

$ws5.Range('A1').Value = 'missingValues, hasLabel, moreFeatures, outliers'
$ws5.Range('B1').Value = 'missingDataFill'
$ws5.Range('C1').Value = 'labelRemove'
$ws5.Range('D1').Value = 'outliersRemove'
$ws5.Range('E1').Value = 'kFold'
$ws5.Range('F1').Value = 'randomForest'
$ws5.Range('G1').Value = 'roc'

$ws5.Range('A2').Value = 'zeroVariance, hasLabel, moreFeatures, outliers'
$ws5.Range('B2').Value = 'zeroVarRemove'
$ws5.Range('C2').Value = 'labelRemove'
$ws5.Range('D2').Value = 'outliersRemove'
$ws5.Range('E2').Value = 'kFold'
$ws5.Range('F2').Value = 'randomForest'
$ws5.Range('G2').Value = 'roc'

$ws5.Range('A3').Value = 'missingValues, zeroVariance, hasLabel, moreFeatures, outliers'
$ws5.Range('B3').Value = 'zeroVarRemove'
$ws5.Range('C3').Value = 'missingDataFill'
$ws5.Range('D3').Value = 'labelRemove'
$ws5.Range('E3').Value = 'outliersRemove'
$ws5.Range('F3').Value = 'kFold'
$ws5.Range('G3').Value = 'randomForest'
$ws5.Range('H3').Value = 'roc'

$ws5.Range('A4').Value = 'hasLabel, moreFeatures, outliers'
$ws5.Range('B4').Value = 'labelRemove'
$ws5.Range('C4').Value = 'outliersRemove'
$ws5.Range('D4').Value = 'kFold'
$ws5.Range('E4').Value = 'randomForest'
$ws5.Range('F4').Value = 'roc'

$ws5.Range('A5').Value = 'missingValues, categorical, hasLabel, moreFeatures, outliers'
$ws5.Range('B5').Value = 'missingValuesRemove'
$ws5.Range('C5').Value = 'labelRemove'
$ws5.Range('D5').Value = 'outliersRemove'
$ws5.Range('E5').Value = 'oneHotEncode'
$ws5.Range('F5').Value = 'kFold'
$ws5.Range('G5').Value = 'randomForest'
$ws5.Range('H5').Value = 'roc'

$ws5.Range('A6').Value = 'categorical, hasLabel, moreFeatures, outliers'
$ws5.Range('B6').Value = 'labelRemove'
$ws5.Range('C6').Value = 'outliersRemove'
$ws5.Range('D6').Value = 'oneHotEncode'
$ws5.Range('E6').Value = 'kFold'
$ws5.Range('F6').Value = 'randomForest'
$ws5.Range('G6').Value = 'roc'

# Re-apply the "bold-ish" black font style (style index 1 in the workbook)
# to the labelRemove/outliersRemove/missingValuesRemove cells, matching the
# styling used throughout the rest of the workbook for these steps.
$styledCells = @('C1','D1','C2','D2','D3','E3','B4','C4','B5','C5','D5','B6','C6')
foreach ($addr in $styledCells) {
    $ws5.Range($addr).Font.Color = 0
}

# This sheet becomes the active one (last created/selected), mirroring the
# workbook's new activeTab.
$ws5.Range('E6').Select()
